$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "bank" keyword/frequency row is inserted right after the header,
# which pushes credit/inflation/interest down by one row and drops the old
# "mortgage"/"recession" rows (only "trade" survives after them).
#
# Rewrite rows 2-6 in place (so each cell keeps its existing per-cell style
# instead of picking up a freshly derived one), then delete the now-unused
# trailing row 7.

$ws.Cells.Item(2, 1).Value = "bank"
$ws.Cells.Item(2, 2).Value = "Frequency"
$ws.Cells.Item(2, 3).Value = 0.0227

$ws.Cells.Item(3, 1).Value = "credit"
$ws.Cells.Item(3, 2).Value = "Frequency"
$ws.Cells.Item(3, 3).Value = 0.0149

$ws.Cells.Item(4, 1).Value = "inflation"
$ws.Cells.Item(4, 2).Value = "Frequency"
$ws.Cells.Item(4, 3).Value = -0.047

$ws.Cells.Item(5, 1).Value = "interest"
$ws.Cells.Item(5, 2).Value = "Frequency"
$ws.Cells.Item(5, 3).Value = -0.0309

$ws.Cells.Item(6, 1).Value = "trade"
$ws.Cells.Item(6, 2).Value = "Frequency"
$ws.Cells.Item(6, 3).Value = -0.1062

# Remove the old trailing row (previously "trade", now redundant).
$ws.Rows.Item(7).Delete()
